$d = $word.ActiveDocument

# --- Edit 1: "/api/lenders" -> "/api/contacts" -------------------------
# The cell holds two runs ("/api/" and "lenders"); only the "lenders"
# run's text should change. A bare Find/Replace would merge the two
# adjacent, identically-formatted runs into one, so we briefly toggle
# Bold on the matched range to keep it a distinct run, then restore it.
$r1 = $d.Content
$r1.Find.Execute("lenders")
$r1.Font.Bold = $true
$r1.Text = "contacts"
$r1.Font.Bold = $false

# --- Edit 2: "Gets all lenders" -> "Gets all contacts to display in lenders list"
# Same situation: "Gets all " and "lenders" are separate runs; keep them
# separate by toggling formatting around the replace.
$r2 = $d.Content
$r2.Find.Execute("lenders")
$r2.Font.Bold = $true
$r2.Text = "contacts to display in lenders list"
$r2.Font.Bold = $false

# --- Edit 3: merge the three runs describing the transactions endpoint --
# "Insert paid information for the " + "specific borrow" + " activity"
# becomes one run with the combined text.
$r3 = $d.Content
$r3.Find.Execute("Insert paid information for the specific borrow activity", $false, $false, $false, $false, $false, $true, 1, $false, "Insert paid information for the specific borrow activity", 2)
